# Updates the "Price" (column D) and "Volume(1h)" (column E) values on the
# cryptos worksheet to reflect the latest scraped figures, matching the
# GitHub Actions scheduled refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a plain-text value into a cell without letting Excel's
# automatic type detection turn numeric-looking strings (e.g. "578.32")
# into real numbers (which would introduce floating point noise and change
# the cell type). We briefly mark the cell as Text, assign the string, then
# restore the cell's original number format/style so no visible formatting
# changes are left behind.
function Set-TextValue {
    param($addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = $origStyle
}

Set-TextValue "D2" "69.893.19"
$ws.Range("E2").Value = "  -0.60%  "

Set-TextValue "D3" "3.578.47"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("E4").Value = "  -0.24%  "

Set-TextValue "D5" "578.32"
$ws.Range("E5").Value = "  -2.02%  "

Set-TextValue "D6" "190.32"
$ws.Range("E6").Value = "  -0.76%  "

Set-TextValue "D7" "0.633"
$ws.Range("E7").Value = "  -2.29%  "

Set-TextValue "D8" "3.575.82"
$ws.Range("E8").Value = "  -0.58%  "

$ws.Range("E9").Value = "  -0.12%  "

Set-TextValue "D10" "0.181"
$ws.Range("E10").Value = "  +0.98%  "

Set-TextValue "D11" "0.662"
$ws.Range("E11").Value = "  -0.03%  "

Set-TextValue "D12" "55.72"
$ws.Range("E12").Value = "  -4.28%  "

$ws.Range("E13").Value = "  +4.15%  "

Set-TextValue "D14" "9.64"
$ws.Range("E14").Value = "  -1.19%  "

Set-TextValue "D15" "4.154.35"
$ws.Range("E15").Value = "  -0.89%  "

Set-TextValue "D16" "19.90"
$ws.Range("E16").Value = "  +2.65%  "

Set-TextValue "D17" "3.575.82"
$ws.Range("E17").Value = "  -1.37%  "

Set-TextValue "D18" "69.827.75"
$ws.Range("E18").Value = "  -0.90%  "

Set-TextValue "D19" "12.67"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("E21").Value = "  -0.89%  "

Set-TextValue "D22" "477.39"
$ws.Range("E22").Value = "  -4.07%  "

Set-TextValue "D23" "19.41"
$ws.Range("E23").Value = "  +12.30%  "

$ws.Range("E24").Value = "  -5.98%  "

Set-TextValue "D25" "96.04"
$ws.Range("E25").Value = "  +5.76%  "

Set-TextValue "D26" "4.38"
$ws.Range("E26").Value = "  -2.57%  "

$ws.Range("E27").Value = "  -3.58%  "

$ws.Range("E28").Value = "  -1.22%  "

$ws.Range("E29").Value = "  -0.33%  "

Set-TextValue "D30" "32.36"
$ws.Range("E30").Value = "  +0.19%  "

Set-TextValue "D31" "7.65"
$ws.Range("E31").Value = "  +1.49%  "

Set-TextValue "D32" "12.24"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("E33").Value = "  +1.46%  "

Set-TextValue "D34" "66.35"
$ws.Range("E34").Value = "  +2.04%  "

Set-TextValue "D35" "581.42"
$ws.Range("E35").Value = "  -6.58%  "

Set-TextValue "D36" "38.95"
$ws.Range("E36").Value = "  +2.27%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  -2.47%  "

Set-TextValue "D39" "0.395"
$ws.Range("E39").Value = "  -2.93%  "

Set-TextValue "D40" "3.25"
$ws.Range("E40").Value = "  +19.95%  "

$ws.Range("E41").Value = "  -6.04%  "

$ws.Range("E42").Value = "  -5.18%  "

Set-TextValue "D43" "2.86"
$ws.Range("E43").Value = "  +6.65%  "

Set-TextValue "D44" "3.227.34"
$ws.Range("E44").Value = "  -2.75%  "

$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("E46").Value = "  -0.77%  "

Set-TextValue "D47" "3.35"
$ws.Range("E47").Value = "  +1.32%  "

Set-TextValue "D48" "9.33"
$ws.Range("E48").Value = "  +2.39%  "

$ws.Range("E49").Value = "  +0.45%  "

$ws.Range("E50").Value = "  -0.21%  "

$ws.Range("E51").Value = "  -5.10%  "
